$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 122.3
$ws.Range("I5").Value = 127.875
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 127.875
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -12.875
$ws.Range("N5").Value = -330
$ws.Range("H132").Value = 1373048.9
$ws.Range("I132").Value = 1056.1666
$ws.Range("J132").Value = 12348990
$ws.Range("K132").Value = 3168.4998
$ws.Range("L132").Value = 37046970
$ws.Range("M132").Value = -638.4998000000001
$ws.Range("N132").Value = -37052030
$ws.Range("H137").Value = 26063270
$ws.Range("I137").Value = 1250.125
$ws.Range("K137").Value = 3750.375
$ws.Range("M137").Value = -1200.375

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 667898.9
$ws.Range("I45").Value = 910294
$ws.Range("J45").Value = 1312.25
$ws.Range("K45").Value = 910294
$ws.Range("L45").Value = 1312.25
$ws.Range("M45").Value = -909917
$ws.Range("N45").Value = -2066.25
$ws.Range("H61").Value = 1933187.4
$ws.Range("I61").Value = 833968.3
$ws.Range("J61").Value = 29413664
$ws.Range("K61").Value = 833968.3
$ws.Range("L61").Value = 29413664
$ws.Range("M61").Value = -833756.3
$ws.Range("N61").Value = -29414088
$ws.Range("H74").Value = 139397760
$ws.Range("I74").Value = 200002130
$ws.Range("J74").Value = 88894120
$ws.Range("K74").Value = 200002130
$ws.Range("L74").Value = 88894120
$ws.Range("M74").Value = -200001256
$ws.Range("N74").Value = -88895868
$ws.Range("H77").Value = 139397760
$ws.Range("I77").Value = 200002130
$ws.Range("J77").Value = 88894120
$ws.Range("K77").Value = 1000010650
$ws.Range("L77").Value = 444470600
$ws.Range("M77").Value = -1000006282
$ws.Range("N77").Value = -444479336
$ws.Range("H136").Value = 1933187.4
$ws.Range("I136").Value = 833968.3
$ws.Range("J136").Value = 29413664
$ws.Range("K136").Value = 2501904.9
$ws.Range("L136").Value = 88240992
$ws.Range("M136").Value = -2499354.9
$ws.Range("N136").Value = -88246092

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1695.5151
$ws.Range("I132").Value = 1067.44
$ws.Range("K132").Value = 3202.32
$ws.Range("M132").Value = -672.3200000000002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4794816
$ws.Range("I5").Value = 6993817.5
$ws.Range("K5").Value = 20981452.5
$ws.Range("M5").Value = -20981340.5
$ws.Range("H33").Value = 72.05556
$ws.Range("I33").Value = 21.25
$ws.Range("J33").Value = 173.66667
$ws.Range("K33").Value = 127.5
$ws.Range("L33").Value = 1042.00002
$ws.Range("M33").Value = 155.5
$ws.Range("N33").Value = -1608.00002
$ws.Range("H97").Value = 1531.0769
$ws.Range("I97").Value = 850.75
$ws.Range("J97").Value = 1833.4445
$ws.Range("K97").Value = 2552.25
$ws.Range("L97").Value = 5500.333500000001
$ws.Range("M97").Value = -2056.25
$ws.Range("N97").Value = -6492.333500000001
$ws.Range("H98").Value = 962.08
$ws.Range("I98").Value = 583.6316
$ws.Range("J98").Value = 2160.5
$ws.Range("K98").Value = 1750.8948
$ws.Range("L98").Value = 6481.5
$ws.Range("M98").Value = -252.8948
$ws.Range("N98").Value = -9477.5
$ws.Range("H132").Value = 2029.1212
$ws.Range("J132").Value = 2027.4348
$ws.Range("L132").Value = 18246.9132
$ws.Range("N132").Value = -23306.9132
$ws.Range("H134").Value = 1446.909
$ws.Range("I134").Value = 872.2222
$ws.Range("J134").Value = 4033
$ws.Range("K134").Value = 2616.6666
$ws.Range("L134").Value = 12099
$ws.Range("M134").Value = 2453.3334
$ws.Range("N134").Value = -22239
$ws.Range("H135").Value = 4794816
$ws.Range("I135").Value = 6993817.5
$ws.Range("K135").Value = 62944357.5
$ws.Range("M135").Value = -62941822.5
$ws.Range("H137").Value = 2734.7058
$ws.Range("J137").Value = 3000
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -19200

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 11111142
$ws.Range("I2").Value = 11111142
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 11111142
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -11111029
$ws.Range("H125").Value = 32312.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 32312.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 32312.5
$ws.Range("N125").Value = -37232.5
$ws.Range("H126").Value = 5405.087
$ws.Range("I126").Value = 7503.533
$ws.Range("J126").Value = 1470.5
$ws.Range("K126").Value = 22510.599
$ws.Range("L126").Value = 4411.5
$ws.Range("M126").Value = -20040.599
$ws.Range("N126").Value = -9351.5
$ws.Range("H127").Value = 34365
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 34365
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 34365
$ws.Range("N127").Value = -44285
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 34254
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 34254
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 34254
$ws.Range("N129").Value = -44254
$ws.Range("H130").Value = 54980
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 54980
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 54980
$ws.Range("N130").Value = -65020
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 7053601
$ws.Range("I132").Value = 7075500.5
$ws.Range("J132").Value = 6994641
$ws.Range("K132").Value = 21226501.5
$ws.Range("L132").Value = 20983923
$ws.Range("M132").Value = -21223971.5
$ws.Range("N132").Value = -20988983
$ws.Range("H133").Value = 51167.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 51167.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 51167.8
$ws.Range("N133").Value = -61287.8
$ws.Range("H134").Value = 18877.166
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 18877.166
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 56631.49800000001
$ws.Range("N134").Value = -61701.49800000001
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 18771.54
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 18771.54
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 56314.62
$ws.Range("N136").Value = -61414.62
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 42999.223
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 42999.223
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 42999.223
$ws.Range("N138").Value = -53279.223
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 77777
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 77777
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 77777
$ws.Range("N140").Value = -88137
$ws.Range("H141").Value = 49500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 49500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 49500
$ws.Range("N141").Value = -59860

Write-Host "Edits applied"